$d = $word.ActiveDocument

$pairs = @(
    @("57×44=2508", "20×85=1700"),
    @("38×24=912", "27×93=2511"),
    @("66×36=2376", "28×12=336"),
    @("33×70=2310", "33×27=891"),
    @("36×35=1260", "64×85=5440"),
    @("40×26=1040", "81×45=3645"),
    @("95×85=8075", "99×55=5445"),
    @("83×99=8217", "55×19=1045"),
    @("54×68=3672", "70×84=5880"),
    @("91×67=6097", "12×82=984"),
    @("30×64=1920", "87×36=3132"),
    @("96×81=7776", "69×41=2829"),
    @("44×11=484", "44×12=528"),
    @("89×57=5073", "41×65=2665"),
    @("20×66=1320", "52×86=4472"),
    @("98×74=7252", "37×60=2220"),
    @("62×51=3162", "45×42=1890"),
    @("13×17=221", "11×37=407"),
    @("77×39=3003", "30×51=1530"),
    @("70×49=3430", "95×72=6840"),
    @("50×87=4350", "14×25=350"),
    @("57×35=1995", "85×69=5865"),
    @("19×17=323", "88×51=4488"),
    @("37×30=1110", "82×59=4838"),
    @("85×37=3145", "75×29=2175")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
